$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15 (notebook 14 / fig_no 15 "Unmet needs for family planning") ---
# J15: "Countries" -> "Percentage"; K15 cell removed (was "Percentage")
$ws.Range("J15").Value = "Percentage"
$ws.Range("K15").ClearContents()
# M15: plot_height 900 -> 640
$ws.Range("M15").Value = 640
# P15: source_link changed to new bit.ly short link
$ws.Range("P15").Value = "http://bit.ly/unmet_need_fp"

# --- Row 16 (notebook 15 / fig_no 14 "Risk of maternal mortality...") ---
# C16: title text update
$ws.Range("C16").Value = "Risk of maternal mortality by Schooling (2005-2012)"
# D16: subtitle text update
$ws.Range("D16").Value = "Maternal mortality ratio compared to the mean number of years in school, females, 2005-2012"
# I16: hovermode x -> closest
$ws.Range("I16").Value = "closest"
# J16, K16: new y_label / x_label values
$ws.Range("J16").Value = "Maternal mortality ratio (deaths per 100,000 live births), 2010"
$ws.Range("K16").Value = "Mean Number of Years in School (Females)"
# O16, P16, Q16: source_label / source_link / source_accessed cleared
$ws.Range("O16").ClearContents()
$ws.Range("P16").ClearContents()
$ws.Range("Q16").ClearContents()

# --- Update sheet view / selection to match author's final cursor position ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("Q16").Select()
